# Update excel processor to handle sheet names and column processing
#
# These are the "section header" rows: the item in column A is a section
# label (e.g. "CI-24HC01713-1S") that was erroneously duplicated into
# columns B, C, D and E as well. Clear the duplicated values in B:E for
# each of those rows, leaving only column A populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sectionRows = @(2,97,194,279,364,417,455,490,525,595,663,730,797,865,932,983,1034,1085,1137,1188,1239,1290,1341,1392,1473,1554,1635,1637,1639,1641,1643)

foreach ($r in $sectionRows) {
    $ws.Range("B$r`:E$r").ClearContents()
}
